# Update numeric profit-calculation cells across several sheets
# (scheduled-runner refresh of market data).
$wb = $excel.ActiveWorkbook

# @@ -1226,25 +1226,25 @@  sheet=ALC row=12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 227.83333
$ws.Range("I12").Value = 241.75
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 241.75
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = -71.75
$ws.Range("N12").Value = -540

# @@ -2258,25 +2258,25 @@  sheet=ALC row=33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 244.65517
$ws.Range("I33").Value = 198.32
$ws.Range("J33").Value = 534.25
$ws.Range("K33").Value = 198.32
$ws.Range("L33").Value = 534.25
$ws.Range("M33").Value = 30.68000000000001
$ws.Range("N33").Value = -992.25

# @@ -3798,19 +3798,22 @@  sheet=ALC row=64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6990
$ws.Range("I64").Value = 6990
$ws.Range("K64").Value = 6990
$ws.Range("M64").Value = -6742

# @@ -3942,19 +3945,22 @@  sheet=ALC row=67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6990
$ws.Range("I67").Value = 6990
$ws.Range("K67").Value = 6990
$ws.Range("M67").Value = -6132

# @@ -4086,22 +4092,22 @@  sheet=ALC row=70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1647.4667
$ws.Range("I70").Value = 1621.5
$ws.Range("K70").Value = 4864.5
$ws.Range("M70").Value = -4594.5

# @@ -4236,22 +4242,22 @@  sheet=ALC row=73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1647.4667
$ws.Range("I73").Value = 1621.5
$ws.Range("K73").Value = 4864.5
$ws.Range("M73").Value = -3928.5

# @@ -5494,22 +5500,22 @@  sheet=ALC row=98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2450.652
$ws.Range("I98").Value = 2524.8572
$ws.Range("K98").Value = 2524.8572
$ws.Range("M98").Value = -1026.8572

# @@ -6700,22 +6706,22 @@  sheet=ALC row=122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2450.652
$ws.Range("I122").Value = 2524.8572
$ws.Range("K122").Value = 7574.571599999999
$ws.Range("M122").Value = -5124.571599999999

# @@ -7049,25 +7055,25 @@  sheet=ALC row=129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 767
$ws.Range("J129").Value = 858.41174
$ws.Range("L129").Value = 2575.23522
$ws.Range("N129").Value = -12575.23522

# @@ -7196,25 +7202,25 @@  sheet=ALC row=132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 9268072
$ws.Range("I132").Value = 15158722
$ws.Range("J132").Value = 11336.143
$ws.Range("K132").Value = 45476166
$ws.Range("L132").Value = 34008.429
$ws.Range("M132").Value = -45473636
$ws.Range("N132").Value = -39068.429

# @@ -7346,25 +7352,25 @@  sheet=ALC row=135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 41668324
$ws.Range("I135").Value = 457.78946
$ws.Range("J135").Value = 200006200
$ws.Range("K135").Value = 4120.105140000001
$ws.Range("L135").Value = 1800055800
$ws.Range("M135").Value = -1585.105140000001
$ws.Range("N135").Value = -1800060870

# @@ -9054,22 +9060,22 @@  sheet=ARM row=28
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6919.25
$ws.Range("I28").Value = 6919.25
$ws.Range("K28").Value = 6919.25
$ws.Range("M28").Value = -6727.25

# @@ -9247,22 +9253,22 @@  sheet=ARM row=32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5286.837
$ws.Range("I32").Value = 4980.3125
$ws.Range("K32").Value = 4980.3125
$ws.Range("M32").Value = -4693.3125

# @@ -10745,25 +10751,25 @@  sheet=ARM row=63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2146.6858
$ws.Range("I63").Value = 2055.1738
$ws.Range("J63").Value = 2322.0833
$ws.Range("K63").Value = 2055.1738
$ws.Range("L63").Value = 2322.0833
$ws.Range("M63").Value = -1369.1738
$ws.Range("N63").Value = -3694.0833

# @@ -10895,25 +10901,25 @@  sheet=ARM row=66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2146.6858
$ws.Range("I66").Value = 2055.1738
$ws.Range("J66").Value = 2322.0833
$ws.Range("K66").Value = 10275.869
$ws.Range("L66").Value = 11610.4165
$ws.Range("M66").Value = -6843.869000000001
$ws.Range("N66").Value = -18474.4165

# @@ -12491,22 +12497,22 @@  sheet=ARM row=99
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 6919.25
$ws.Range("I99").Value = 6919.25
$ws.Range("K99").Value = 6919.25
$ws.Range("M99").Value = -3924.25

# @@ -18543,22 +18549,22 @@  sheet=BSM row=82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17013.334
$ws.Range("I82").Value = 10020
$ws.Range("K82").Value = 10020
$ws.Range("M82").Value = -9637

# @@ -18696,22 +18702,22 @@  sheet=BSM row=85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 17013.334
$ws.Range("I85").Value = 10020
$ws.Range("K85").Value = 10020
$ws.Range("M85").Value = -8694

# @@ -18748,25 +18754,25 @@  sheet=BSM row=86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2919.889
$ws.Range("I86").Value = 2967.2
$ws.Range("J86").Value = 2683.3333
$ws.Range("K86").Value = 2967.2
$ws.Range("L86").Value = 2683.3333
$ws.Range("M86").Value = -1844.2
$ws.Range("N86").Value = -4929.3333

# @@ -18898,25 +18904,25 @@  sheet=BSM row=89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2919.889
$ws.Range("I89").Value = 2967.2
$ws.Range("J89").Value = 2683.3333
$ws.Range("K89").Value = 14836
$ws.Range("L89").Value = 13416.6665
$ws.Range("M89").Value = -9220
$ws.Range("N89").Value = -24648.6665

# @@ -21091,25 +21097,25 @@  sheet=BSM row=134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7795.294
$ws.Range("I134").Value = 1229.3636
$ws.Range("J134").Value = 19832.834
$ws.Range("K134").Value = 3688.0908
$ws.Range("L134").Value = 59498.50199999999
$ws.Range("M134").Value = -1153.0908
$ws.Range("N134").Value = -64568.50199999999

# @@ -22989,25 +22995,25 @@  sheet=CRP row=31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1067.34
$ws.Range("I31").Value = 1026.0975
$ws.Range("J31").Value = 1255.2222
$ws.Range("K31").Value = 1026.0975
$ws.Range("L31").Value = 1255.2222
$ws.Range("M31").Value = -731.0975000000001
$ws.Range("N31").Value = -1845.2222

# @@ -23136,25 +23142,25 @@  sheet=CRP row=34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1067.34
$ws.Range("I34").Value = 1026.0975
$ws.Range("J34").Value = 1255.2222
$ws.Range("K34").Value = 1026.0975
$ws.Range("L34").Value = 1255.2222
$ws.Range("M34").Value = -824.0975000000001
$ws.Range("N34").Value = -1659.2222

# @@ -26500,22 +26506,22 @@  sheet=CRP row=104
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 28500
$ws.Range("J104").Value = 28500
$ws.Range("L104").Value = 28500
$ws.Range("N104").Value = -33742

# @@ -27376,25 +27382,25 @@  sheet=CRP row=122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 817.3684
$ws.Range("I122").Value = 833.8461
$ws.Range("J122").Value = 781.6667
$ws.Range("K122").Value = 2501.5383
$ws.Range("L122").Value = 2345.0001
$ws.Range("M122").Value = -51.53830000000016
$ws.Range("N122").Value = -7245.0001

# @@ -27872,25 +27878,25 @@  sheet=CRP row=132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 8808.866
$ws.Range("I132").Value = 15774.571
$ws.Range("J132").Value = 2713.875
$ws.Range("K132").Value = 47323.713
$ws.Range("L132").Value = 8141.625
$ws.Range("M132").Value = -44793.713
$ws.Range("N132").Value = -13201.625

# @@ -27973,25 +27979,25 @@  sheet=CRP row=134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 25002706
$ws.Range("I134").Value = 3253.923
$ws.Range("J134").Value = 71430260
$ws.Range("K134").Value = 9761.769
$ws.Range("L134").Value = 214290780
$ws.Range("M134").Value = -7226.769
$ws.Range("N134").Value = -214295850

# @@ -35029,25 +35035,25 @@  sheet=CUL row=131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16952166
$ws.Range("I131").Value = 111111416
$ws.Range("J131").Value = 3500.72
$ws.Range("K131").Value = 333334248
$ws.Range("L131").Value = 10502.16
$ws.Range("M131").Value = -333329208
$ws.Range("N131").Value = -20582.16

# @@ -48364,22 +48370,22 @@  sheet=LTW row=122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 41690800
$ws.Range("I122").Value = 62525624
$ws.Range("K122").Value = 187576872
$ws.Range("M122").Value = -187574422

# @@ -48854,25 +48860,25 @@  sheet=LTW row=132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 113140.6
$ws.Range("I132").Value = 37067.668
$ws.Range("J132").Value = 145743.28
$ws.Range("K132").Value = 111203.004
$ws.Range("L132").Value = 437229.84
$ws.Range("M132").Value = -108673.004
$ws.Range("N132").Value = -442289.84

# @@ -52306,25 +52312,25 @@  sheet=WVR row=62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 71439860
$ws.Range("I62").Value = 125008750
$ws.Range("J62").Value = 14667.333
$ws.Range("K62").Value = 125008750
$ws.Range("L62").Value = 14667.333
$ws.Range("M62").Value = -125008126
$ws.Range("N62").Value = -15915.333

# @@ -52456,25 +52462,25 @@  sheet=WVR row=65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 71439860
$ws.Range("I65").Value = 125008750
$ws.Range("J65").Value = 14667.333
$ws.Range("K65").Value = 625043750
$ws.Range("L65").Value = 73336.66500000001
$ws.Range("M65").Value = -625040630
$ws.Range("N65").Value = -79576.66500000001

# @@ -55409,25 +55415,25 @@  sheet=WVR row=126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 58831636
$ws.Range("I126").Value = 71430840
$ws.Range("J126").Value = 35335
$ws.Range("K126").Value = 214292520
$ws.Range("L126").Value = 106005
$ws.Range("M126").Value = -214290050
$ws.Range("N126").Value = -110945

# @@ -55703,25 +55709,25 @@  sheet=WVR row=132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5278
$ws.Range("I132").Value = 9001.333000000001
$ws.Range("J132").Value = 3416.3333
$ws.Range("K132").Value = 27003.999
$ws.Range("L132").Value = 10248.9999
$ws.Range("M132").Value = -24473.999
$ws.Range("N132").Value = -15308.9999

# @@ -55899,25 +55905,25 @@  sheet=WVR row=136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 966.3823
$ws.Range("I136").Value = 861.125
$ws.Range("J136").Value = 1219
$ws.Range("K136").Value = 2583.375
$ws.Range("L136").Value = 3657
$ws.Range("M136").Value = -33.375
$ws.Range("N136").Value = -8757
